$d = $word.ActiveDocument

# --- Create the three new character styles -------------------------------
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "campaign period" run ------------------------
$periodText = "Perioadele campaniei din 2022 pentru Constelația Leului: 14-23 aprilie, 14-23 mai"
$range = $d.Content
$range.Start = 0
$range.End = $d.Content.End
while ($range.Find.Execute($periodText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $range.Style = "GaNStyle"
    $range.Start = $range.End
    $range.End = $d.Content.End
}

# --- Apply GaNParagraph to the campaign description paragraph -------------
$paragraphText = "Prin această activitate participați în cadrul unei campanii globale de observare și consemnare a celor mai slabe stele vizibile ca metodă de măsurare a poluării luminoase dintr-un anumit loc. Localizând și observând  Constelația Leului pe cerul nopții și comparând-o cu diagramele stelare, oamenii din întreaga lume vor putea afla în ce măsură iluminatul nocturn din comunitatea lor contribuie la poluarea luminoasă. Contribuțiile dumneavoastră la baza de date online vor facilita o documentare globală privind cerul nocturn observabil."
$range2 = $d.Content
$range2.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$range2.Style = "GaNParagraph"

# --- Apply GaNLinks to the credit/link run ---------------------------------
$linkText = "de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$range3 = $d.Content
$range3.Find.Execute($linkText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$range3.Style = "GaNLinks"
